# Update cryptocurrency price/volume data in the worksheet.
# Source data refresh as of Thu Jul 27 16:09:31 UTC 2023 (GitHub Actions).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New cell values scraped for this run (Coin/Link/Price/Volume columns).
$cellValues = @{
    "D2" = "29.320.31"
    "E2" = "  -0.05%  "
    "D3" = "1.870.00"
    "E3" = "  +0.44%  "
    "D4" = "1.001"
    "E4" = "  +0.08%  "
    "D5" = "0.7213"
    "E5" = "  +2.38%  "
    "D6" = "240.84"
    "E6" = "  +1.07%  "
    "D7" = "1.001"
    "E7" = "  +0.04%  "
    "D8" = "0.07880"
    "E8" = "  -0.04%  "
    "D9" = "0.3091"
    "E9" = "  +1.26%  "
    "D10" = "25.25"
    "E10" = "  +1.58%  "
    "D11" = "0.08261"
    "E11" = "  +0.95%  "
    "D12" = "0.7221"
    "E12" = "  +0.57%  "
    "D13" = "1.858.13"
    "E13" = "  -6.62%  "
    "D14" = "5.239"
    "E14" = "  +0.20%  "
    "D15" = "90.71"
    "E15" = "  +1.18%  "
    "D16" = "29.328.88"
    "E16" = "  -0.03%  "
    "D17" = "5.846"
    "E17" = "  -0.11%  "
    "D18" = "243.60"
    "E18" = "  +2.10%  "
    "E19" = "  +0.17%  "
    "D20" = "13.22"
    "E20" = "  -0.32%  "
    "D21" = "2.112.07"
    "E21" = "  -6.49%  "
    "D22" = "1.001"
    "E22" = "  -0.01%  "
    "E23" = "  +4.92%  "
    "D24" = "1.002"
    "E24" = "  +0.07%  "
    "D25" = "0.1612"
    "E25" = "  +12.29%  "
    "D26" = "162.50"
    "E26" = "  -0.30%  "
    "D27" = "8.967"
    "E27" = "  +0.52%  "
    "D28" = "18.25"
    "E28" = "  +0.81%  "
    "E29" = "  -2.01%  "
    "E30" = "  +1.22%  "
    "D31" = "4.374"
    "E31" = "  +1.03%  "
    "D32" = "4.112"
    "E32" = "  +1.24%  "
    "D33" = "0.05207"
    "E33" = "  +0.22%  "
    "D34" = "1.938"
    "E34" = "  +1.71%  "
    "D35" = "1.186"
    "E35" = "  +0.59%  "
    "D36" = "0.7277"
    "E36" = "  +2.10%  "
    "D37" = "2.684"
    "E37" = "  +0.30%  "
    "E38" = "  +0.35%  "
    "D39" = "2.700"
    "E39" = "  +0.30%  "
    "D40" = "1.174.09"
    "E40" = "  -0.21%  "
    "D41" = "0.9031"
    "E41" = "  -2.05%  "
    "D42" = "6.123"
    "E42" = "  +1.47%  "
    "D43" = "72.44"
    "E43" = "  +0.69%  "
    "D44" = "1.001"
    "E44" = "  +0.02%  "
    "D45" = "101.64"
    "E45" = "  -0.50%  "
    "D46" = "0.5284"
    "E46" = "  -1.18%  "
    "D47" = "2.010.95"
    "E47" = "  -6.38%  "
    "D48" = "1.782"
    "E48" = "  +1.09%  "
    "B49" = "BabyDogeCoin"
    "C49" = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
    "D49" = "0.00000000121"
    "E49" = "  +1.88%  "
    "B50" = "SynthetixNetwork"
    "C50" = "https://coinranking.com/coin/sgxZRXbK0FDc+synthetixnetwork-snx"
    "D50" = "2.893"
    "E50" = "  +5.50%  "
    "B51" = "EnergySwap"
    "C51" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "D51" = "9.254"
    "E51" = "  +0.57%  "
}

# Several Price values are purely numeric-looking strings (e.g. "1.001").
# Force those specific cells to Text format first so Excel keeps them as
# literal strings instead of silently converting them to numbers.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

foreach ($addr in $cellValues.Keys) {
    $ws.Range($addr).Value = $cellValues[$addr]
}
